# Weekly refresh of fruit/vegetable prices: rows 2-8 are rotated up by
# three (data for rows 5-8 moves to rows 2-5, and data for rows 2-4 moves
# down to rows 6-8), reflecting the latest market report dates.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot current values for rows 2-8 (columns D, I, J, K, L, M, N, P, Q)
$cols = @("D", "I", "J", "K", "L", "M", "N", "P", "Q")
$rows = 2..8

$snapshot = @{}
foreach ($r in $rows) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowData
}

# New row order: after[r] = before[r+3], wrapping within 2..8
$mapping = @{ 2 = 5; 3 = 6; 4 = 7; 5 = 8; 6 = 2; 7 = 3; 8 = 4 }

foreach ($r in $rows) {
    $src = $mapping[$r]
    $data = $snapshot[$src]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value = $data[$c]
    }
}
